# Commit: "The model for 3rd iteration for CCAC CBWG"
#
# The underlying edit overwrites the raw Adjustment-Factor (AF) data held in
# the "BY-RSD-SH_AF" (space heating) and "BY-RSD-WH_AF" (water heating)
# sheets (columns L and M). The "AF" sheet's D column consumes these via
# AVERAGEIFS(...) / direct cell references, so simply correcting the source
# data and letting Excel recalculate reproduces all of the downstream
# cached formula results in the "AF" sheet as well.

$wb = $excel.ActiveWorkbook

# --- BY-RSD-SH_AF (space heating) : row -> new value for columns L and M ---
$wsSH = $wb.Worksheets.Item("BY-RSD-SH_AF")

$shChanges = @{
    2  = 0.12
    3  = 0.12
    6  = 0.12
    8  = 0.12
    10 = 0.12
    12 = 0.12
    13 = 0.12
    14 = 0.12
    18 = 0.12
    20 = 0.12
    25 = 0.12
    26 = 0.0697946718860816
    30 = 0.0697946718860816
    32 = 0.12
}

foreach ($r in $shChanges.Keys) {
    $val = $shChanges[$r]
    $wsSH.Cells.Item($r, 12).Value = $val   # column L
    $wsSH.Cells.Item($r, 13).Value = $val   # column M
}

# --- BY-RSD-WH_AF (water heating) : row -> new value for columns L and M ---
$wsWH = $wb.Worksheets.Item("BY-RSD-WH_AF")

$whChanges = @{
    2  = 0.12
    3  = 0.12
    6  = 0.12
    8  = 0.12
    9  = 0.12
    10 = 0.12
    11 = 0.12
    12 = 0.12
    14 = 0.12
    15 = 0.12
    19 = 0.12
    21 = 0.12
    23 = 0.12
    24 = 0.12
    27 = 0.12
    28 = 0.12
    29 = 0.12
    30 = 0.12
    31 = 0.12
    32 = 0.12
    34 = 0.12
    36 = 0.12
    37 = 0.12
    40 = 0.12
}

foreach ($r in $whChanges.Keys) {
    $val = $whChanges[$r]
    $wsWH.Cells.Item($r, 12).Value = $val   # column L
    $wsWH.Cells.Item($r, 13).Value = $val   # column M
}

# Force a full recalculation so the "AF" sheet's cached formula results
# (D6, D7, ..., D49, D50, ..., D92, ... which reference the sheets above via
# AVERAGEIFS or direct cell references) pick up the new source values.
$excel.CalculateFullRebuild()
